$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates: force text so Excel does not coerce numeric-looking strings ---
$dCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($ref in $dCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = "27.652.42"
$ws.Range("D3").Value = "1.879.28"
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "331.25"
$ws.Range("D7").Value = "0.4719"
$ws.Range("D8").Value = "0.3962"
$ws.Range("D9").Value = "48.21"
$ws.Range("D10").Value = "0.08034"
$ws.Range("D11").Value = "1.025"
$ws.Range("D13").Value = "1.876.03"
$ws.Range("D14").Value = "5.963"
$ws.Range("D15").Value = "7.166"
$ws.Range("D16").Value = "1.005"
$ws.Range("D17").Value = "87.01"
$ws.Range("D18").Value = "0.00001042"
$ws.Range("D19").Value = "0.06625"
$ws.Range("D20").Value = "17.16"
$ws.Range("D21").Value = "1.003"
$ws.Range("D22").Value = "27.662.72"
$ws.Range("D23").Value = "5.494"
$ws.Range("D24").Value = "11.00"
$ws.Range("D25").Value = "2.291"
$ws.Range("D26").Value = "2.116.57"
$ws.Range("D27").Value = "156.17"
$ws.Range("D28").Value = "20.25"
$ws.Range("D29").Value = "2.088"
$ws.Range("D30").Value = "5.571"
$ws.Range("D31").Value = "122.24"
$ws.Range("D32").Value = "0.09550"
$ws.Range("D33").Value = "0.9644"
$ws.Range("D34").Value = "1.456"
$ws.Range("D35").Value = "3.630"
$ws.Range("D36").Value = "5.301"
$ws.Range("D37").Value = "0.06111"
$ws.Range("D38").Value = "0.02251"
$ws.Range("D39").Value = "1.227"
$ws.Range("D40").Value = "8.106"
$ws.Range("D41").Value = "0.6004"
$ws.Range("D42").Value = "1.002"
$ws.Range("D43").Value = "0.1898"
$ws.Range("D44").Value = "10.26"
$ws.Range("D45").Value = "0.5707"
$ws.Range("D46").Value = "1.251"
$ws.Range("D47").Value = "12.15"
$ws.Range("D48").Value = "3.407"
$ws.Range("D49").Value = "1.933"
$ws.Range("D50").Value = "0.06820"
$ws.Range("D51").Value = "111.04"

foreach ($ref in $dCells) { $ws.Range($ref).Style = "Normal" }

# --- Columns B, C, E updates (plain text; no coercion risk) ---
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("E40").Value = "  -5.21%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  +0.90%  "
